# Update the date line
$d = $word.ActiveDocument

$wdReplaceAll = 2

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, $wdReplaceAll)
}

Replace-Text "2023-07-23 Sunday" "2023-07-24 Monday"

Replace-Text "85÷9=" "14÷7="
Replace-Text "70÷8=" "21÷7="
Replace-Text "49÷2=" "61÷9="
Replace-Text "58÷7=" "66÷4="
Replace-Text "28÷8=" "21÷6="
Replace-Text "82÷2=" "10÷2="
Replace-Text "84÷8=" "83÷8="
Replace-Text "46÷6=" "63÷5="
Replace-Text "68÷4=" "92÷8="
Replace-Text "65÷9=" "43÷7="
Replace-Text "86÷4=" "30÷5="
Replace-Text "59÷7=" "73÷7="
Replace-Text "93÷4=" "76÷5="
Replace-Text "22÷2=" "81÷4="
Replace-Text "58÷2=" "29÷4="
Replace-Text "24÷5=" "30÷3="
Replace-Text "85÷3=" "94÷7="
Replace-Text "22÷8=" "60÷3="
Replace-Text "94÷3=" "80÷5="
Replace-Text "35÷6=" "33÷5="
Replace-Text "10÷6=" "10÷2="
Replace-Text "76÷4=" "19÷9="
Replace-Text "42÷2=" "24÷9="
Replace-Text "79÷9=" "94÷4="
Replace-Text "45÷3=" "57÷5="
